# This script applies the SectorGroup.xlsx edit: the columns
# "codeforiati:category-name" (D) and "codeforiati:group-code" (G) swap places
# with each other (both the header label in row 1 and every data value),
# while columns E ("codeforiati:group-name") and F ("codeforiati:category-code")
# stay exactly where they are.
#
# Net effect per row: new D = old G, new G = old D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Determine the extent of the used data (dimension is A1:G235 in this workbook,
# but compute it dynamically to be safe).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Read column D and column G (rows 1..lastRow) in one shot each, swap them,
# and write back in one shot each -- much faster than per-cell COM calls.
$rangeD = $ws.Range("D1:D$lastRow")
$rangeG = $ws.Range("G1:G$lastRow")

$valuesD = $rangeD.Value2
$valuesG = $rangeG.Value2

# Force text format on both columns so that purely-numeric-looking codes
# (e.g. "110", "998") remain stored as text, exactly like the rest of the
# codeforiati code/group columns in this sheet, rather than being
# reinterpreted as numbers.
$rangeD.NumberFormat = "@"
$rangeG.NumberFormat = "@"

$rangeD.Value2 = $valuesG
$rangeG.Value2 = $valuesD

# Restore the original (default/General) number format so no new cell style
# is introduced -- only the values of D and G should change.
$rangeD.NumberFormat = "General"
$rangeG.NumberFormat = "General"
